$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New weekly WALCL observations appended to the Data sheet (rows 104-109),
# matching the style (date format / border / alignment) already used by
# the existing data rows (e.g. row 103).
$newRows = @(
    @{ Row = 104; Date = 45189; Value = 8024.09 },
    @{ Row = 105; Date = 45196; Value = 8002.064 },
    @{ Row = 106; Date = 45203; Value = 7955.782 },
    @{ Row = 107; Date = 45210; Value = 7952.054 },
    @{ Row = 108; Date = 45217; Value = 7933.162 },
    @{ Row = 109; Date = 45224; Value = 7907.83 }
)

foreach ($r in $newRows) {
    $srcRow = $r.Row - 1
    $ws.Range("A${srcRow}:B${srcRow}").Copy()
    $ws.Range("A$($r.Row):B$($r.Row)").Insert()
    $ws.Range("A$($r.Row)").Borders.LineStyle = 1
    $ws.Range("A$($r.Row)").Value = $r.Date
    $ws.Range("B$($r.Row)").Value = $r.Value
}

# Refresh FRED series metadata on the SeriesInfo sheet.
$info = $wb.Worksheets.Item("SeriesInfo")
$info.Range("B3").Value = "2023-10-27"
$info.Range("B4").Value = "2023-10-27"
$info.Range("B7").Value = "2023-10-25"
$info.Range("B14").Value = "2023-10-26 15:33:02-05"
$info.Range("B15").Value = 93
